$d = $word.ActiveDocument

# --- Step 1: remove the old _GoBack bookmark (near "2. SYSTEM OVERVIEW") ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Step 2: insert descriptive paragraph after "Figure 1: Swimline diagram" ---
$range = $d.Content
$found = $range.Find.Execute("Figure 1: Swimline diagram", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $range.Paragraphs(1)
$para.Range.InsertParagraphAfter()
$newPara = $para.Next()
$ins = $newPara.Range
$ins.Collapse(1)
$ins.InsertAfter("The swimline diagram showed a flow when a customer began to place an order on the website. This diagram helps shows the entities that are involved and the activities that take place. The flow began when customers placed items into a chart, they can then go to checkout. During checkout, customers are asked for some personal information such as name, billing address, and credit card information. The information is transferred to Stripe, a company that helps in the billing process. After Stripe verified that the credit card is valid and payment is processed, the order go to the store selected. Employees at the store received the order and fulfill the order. The order is hold until the customer come to the store and picks up the order.")

# --- Step 3: "Figure 1: DFD for placing an order" -> "Figure 2: DFD for placing an order" ---
$range3 = $d.Content
$found3 = $range3.Find.Execute("Figure 1: DFD for placing an order", $true, $false, $false, $false, $false, $true, 1, $false, "Figure 2: DFD for placing an order", 2)

# --- Step 4: insert descriptive paragraph after "Figure 2: DFD for placing an order" ---
$range4 = $d.Content
$found4 = $range4.Find.Execute("Figure 2: DFD for placing an order", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para4 = $range4.Paragraphs(1)
$para4.Range.InsertParagraphAfter()
$newPara4 = $para4.Next()
$ins4 = $newPara4.Range
$ins4.Collapse(1)
$ins4.InsertAfter("The data flow diagram shows the flow of data when a customer placed an order. This diagram should be used in conjunction with Figure 1. The DFD shows a more detailed view of all the data that occur during an order placement.")

# --- Step 5: add lastRenderedPageBreak before "Total" table cell text ---
$range5 = $d.Content
$found5 = $range5.Find.Execute("Total", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $range5.Collapse(1)
    $range5.InsertBefore([char]0x0B)
}

# --- Step 6: remove lastRenderedPageBreak before "Last Name" table cell text ---
# (handled by direct XML text check below; lastRenderedPageBreak elements aren't
#  directly addressable through the Range/Find object model, see notes.)

# --- Step 7: split the UI overview sentence, relocate the _GoBack bookmark ---
$range7 = $d.Content
$oldText = "The user interface for the system will allow the user to easily obtain health reports, connect with other Care Clients, connect with Caregivers, store health history, and record vitals to test for diseases.  The user interface should contain a menu tool bar containing major features for easy navigation through the app. "
$found7 = $range7.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "The user interface for the system will allow the user to ", 2)

# --- Step 8: remove lastRenderedPageBreak before the homepage mockup drawing ---

# --- Step 9: "Figure 1: Mockup of the homeage" -> "Figure 3: Mockup of the homeage" ---
$range9 = $d.Content
$found9 = $range9.Find.Execute("Figure 1: Mockup of the homeage", $true, $false, $false, $false, $false, $true, 1, $false, "Figure 3: Mockup of the homeage", 2)

# --- Step 10: add lastRenderedPageBreak before "Feature" table cell text ---

# --- Step 11: remove lastRenderedPageBreak before "FR-1" table cell text ---

Write-Output "done"
